$wb = $excel.ActiveWorkbook

# The workbook reports status "In Translation" with timestamps right before
# handoff; the report is being regenerated to reflect "Ready for handoff"
# with refreshed handoff timestamps, and the Status/Datetime columns are
# widened slightly to fit the new (longer) status text.
#
# NOTE: Excel's COM ColumnWidth setter always quantizes to whole on-screen
# pixels, so the widest value reachable is used (nearest achievable width).

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-24 22:39:50"
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-24 22:39:46"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-24 22:39:50"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
